$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F7:I7 values (Xception row) reflecting base model trainable change
$ws.Range("F7").Value = 0.9109
$ws.Range("G7").Value = 0.92132353782653797
$ws.Range("H7").Value = 0.371
$ws.Range("I7").Value = 0.29586303234100297

# Update view: zoom normal and selection
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("E7").Select()
